$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every Price (D) / Volume(1h) (E) cell in the source sheet is stored as
# plain text (inline string), even when the text looks like a plain
# number (e.g. "348.08"). Pre-format every cell we are about to touch as
# Text ("@") *before* writing its value so Excel doesn't reinterpret the
# string as a number (which would, among other things, silently drop
# trailing zeros like "1.010" -> 1.01).
$textCells = "D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8", `
  "D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","E15", `
  "D16","E16","D17","E17","E18","D19","D20","E20", `
  "B21","C21","D21","E21","B22","C22","D22","E22", `
  "D23","E23","D24","E24","D25","E25","D26","E26","E27","D28","E28","D29","E29", `
  "D30","E30","D31","E31","D32","E32","E33","D34","E34", `
  "B35","C35","D35","E35","B36","C36","D36","E36", `
  "D37","E37","D38","E38","D39","E39", `
  "B40","C40","D40","E40","B41","C41","D41","E41", `
  "E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48", `
  "D49","E49","D50","E50","D51","E51"

foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 21 <-> Row 22 swap: Uniswap/Dai change places (with new values)
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.007"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "6.281"
$ws.Range("E22").Value = "  +1.14%  "

# Row 35 <-> Row 36 swap: HuobiToken/InternetComputer(DFINITY) change places (with new values)
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "6.564"
$ws.Range("E35").Value = "  +6.12%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "3.969"
$ws.Range("E36").Value = "  +0.31%  "

# Row 40 <-> Row 41 swap: TheSandbox/Aptos change places (with new values)
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "12.70"
$ws.Range("E40").Value = "  +1.79%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.7026"
$ws.Range("E41").Value = "  +1.90%  "

# Price (D) / Volume(1h) (E) refresh for every other row
$ws.Range("D2").Value = "29.953.09"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").Value = "2.115.59"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.64%  "

$ws.Range("D5").Value = "347.51"
$ws.Range("E5").Value = "  +0.31%  "

$ws.Range("E6").Value = "  +0.50%  "

$ws.Range("D7").Value = "0.5200"
$ws.Range("E7").Value = "  +0.95%  "

$ws.Range("D8").Value = "0.4449"
$ws.Range("E8").Value = "  +0.75%  "

$ws.Range("D9").Value = "54.37"
$ws.Range("E9").Value = "  +4.63%  "

$ws.Range("D10").Value = "0.09349"
$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("D11").Value = "1.178"
$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").Value = "25.15"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").Value = "2.157.82"
$ws.Range("E13").Value = "  +2.78%  "

$ws.Range("D14").Value = "8.455"
$ws.Range("E14").Value = "  +3.85%  "

$ws.Range("E15").Value = "  +2.06%  "

$ws.Range("D16").Value = "102.25"
$ws.Range("E16").Value = "  +2.68%  "

$ws.Range("D17").Value = "0.00001162"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("E18").Value = "  +0.39%  "

$ws.Range("D19").Value = "21.53"

$ws.Range("D20").Value = "0.06687"
$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("D23").Value = "29.991.94"
$ws.Range("E23").Value = "  -0.30%  "

$ws.Range("D24").Value = "12.74"
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("D25").Value = "2.330"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").Value = "2.347.90"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").Value = "  +0.74%  "

$ws.Range("D28").Value = "2.546"
$ws.Range("E28").Value = "  -0.20%  "

$ws.Range("D29").Value = "162.57"
$ws.Range("E29").Value = "  +0.21%  "

$ws.Range("D30").Value = "133.89"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").Value = "1.153"
$ws.Range("E31").Value = "  -1.22%  "

$ws.Range("D32").Value = "1.778"
$ws.Range("E32").Value = "  +8.73%  "

$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").Value = "6.237"
$ws.Range("E34").Value = "  +0.29%  "

$ws.Range("D37").Value = "10.80"
$ws.Range("E37").Value = "  +7.33%  "

$ws.Range("D38").Value = "0.02603"
$ws.Range("E38").Value = "  +1.66%  "

$ws.Range("D39").Value = "0.06853"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").Value = "1.333"
$ws.Range("E43").Value = "  +1.70%  "

$ws.Range("D44").Value = "0.6840"
$ws.Range("E44").Value = "  +3.27%  "

$ws.Range("D45").Value = "14.52"
$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("D46").Value = "2.352"
$ws.Range("E46").Value = "  +3.40%  "

$ws.Range("D47").Value = "1.007"
$ws.Range("E47").Value = "  +0.49%  "

$ws.Range("D48").Value = "3.633"
$ws.Range("E48").Value = "  +0.08%  "

$ws.Range("D49").Value = "0.00000000357"
$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").Value = "1.247"
$ws.Range("E50").Value = "  +7.13%  "

$ws.Range("D51").Value = "1.222"
$ws.Range("E51").Value = "  +0.39%  "
